$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column I ("roboticS1Prep") held text "No" for every data row (2-27).
# Convert those cells to a real Boolean FALSE value, displayed via a
# custom "TRUE/FALSE" number format, as part of cleaning up the database.
$range = $ws.Range("I2:I27")
$range.NumberFormat = '"TRUE";"TRUE";"FALSE"'
$range.Value = $false

# Reflect the new active selection used while making this edit.
$ws.Range("I2:I27").Select()
